$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new data row at 21 (pushes the old total row 21 -> 22
#     and the old footer row 22 -> 23), matching the existing item-row
#     pattern (rows 7-20) for styling/borders.
$ws.Rows("21").Insert()

# Copy the formatting (styles/borders/fills) from the row above (row 20)
# onto the freshly inserted row 21 so it matches the other item rows.
$ws.Range("A20:Q20").Copy()
$ws.Range("A21:Q21").PasteSpecial(-4122)
$ws.Rows("21").RowHeight = 25.5

# --- Populate the new item row (#15) ---
$ws.Range("A21").Value = 15
$ws.Range("C21").Value = "فازلين بيور صغير "
$ws.Range("H21").Value = "5:0"

# L21 and P21 keep their existing numeric display formats but hold
# text values (matching the pattern already used by rows 7-20), so we
# briefly force a text format, set the value, then restore the number
# format to keep the same style id.
$lFmt = $ws.Range("L21").NumberFormat
$ws.Range("L21").NumberFormat = "@"
$ws.Range("L21").Value = "0"
$ws.Range("L21").NumberFormat = $lFmt

$ws.Range("N21").Value = "10.00"

$pFmt = $ws.Range("P21").NumberFormat
$ws.Range("P21").NumberFormat = "@"
$ws.Range("P21").Value = "10.0000"
$ws.Range("P21").NumberFormat = $pFmt

$ws.Range("Q21").Value = "1:0"

# --- Merge the new row's cells the same way as the other item rows ---
$ws.Range("A21:B21").Merge()
$ws.Range("C21:G21").Merge()
$ws.Range("H21:K21").Merge()
$ws.Range("L21:M21").Merge()
$ws.Range("N21:O21").Merge()

# --- Update the running total (old row 21, now row 22) ---
$ws.Range("P22").Value = 673.91999999999996

# --- Update the generated timestamp in the footer (old row 22, now row 23) ---
$ws.Range("A23").Value = "Wednesday, 27 August, 2025 11:00 AM"
